$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format first so Excel does not
# auto-convert numeric-looking strings (e.g. "61.01", "0.0560") into
# floating point numbers and lose formatting / trailing zeros.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '60.907.15'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '2.594.46'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '523.21'
$ws.Range("E5").Value = '  +3.20%  '
$ws.Range("D6").Value = '154.72'
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").Value = '  +2.03%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D13").Value = '3.050.25'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '60.922.28'
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").Value = '21.66'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("D17").Value = '2.599.90'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '353.29'
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").Value = '10.59'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = '6.23'
$ws.Range("E21").Value = '  +1.60%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '61.01'
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("E24").Value = '  +1.74%  '
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = '2.713.43'
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '0.0₃0846'
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = '6.30'
$ws.Range("E31").Value = '  +10.70%  '
$ws.Range("D32").Value = '19.35'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").Value = '1.60'
$ws.Range("E33").Value = '  +3.09%  '
$ws.Range("D34").Value = '148.02'
$ws.Range("E34").Value = '  -3.79%  '
$ws.Range("D35").Value = '4.19'
$ws.Range("E35").Value = '  +5.04%  '
$ws.Range("D36").Value = '0.934'
$ws.Range("E36").Value = '  +8.62%  '
$ws.Range("E37").Value = '  +1.06%  '
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("E41").Value = '  +1.62%  '
$ws.Range("D42").Value = '287.57'
$ws.Range("E42").Value = '  -3.17%  '
$ws.Range("E43").Value = '  +1.47%  '
$ws.Range("E44").Value = '  +1.35%  '
$ws.Range("D45").Value = '0.0560'
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = '19.67'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '4.89'
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0238'
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").Value = '19.10'
$ws.Range("E51").Value = '  +8.64%  '


# Restore the default (unformatted) style on the touched range so the
# saved file does not carry a leftover explicit "Text" number format on
# these cells (matches the original workbook's un-styled inline strings).
$dataRange.Style = "Normal"
